$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell A2 currently holds the number 1; change it to the text "teste"
$ws.Range("A2").Value = "teste"
